$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the split "For some time p" / "eople have dreamed..." runs (which
#    were separated by the _GoBack bookmark) into a single run, and drop the
#    now-redundant bookmark from this spot (it gets re-created further down).
# ---------------------------------------------------------------------------
$mergedIntro = "For some time people have dreamed of a car that drives itself safely on the highway and city streets. While driving can often be a fun experience, Americans spend about 50 billion hours per year behind the wheel and much of that is not productive or very pleasant."
$introPara = $d.Paragraphs.Item(8)
$found = $introPara.Range.Find.Execute($mergedIntro, $true, $false, $false, $false, $false, $true, 1, $false, $mergedIntro, 2)

# ---------------------------------------------------------------------------
# 2. Replace the single trailing empty paragraph (right after "- Maybe
#    briefly list...") with: 7 empty paragraphs, Drew's whole new outline
#    block, and a paragraph holding just the relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$blankPara = $d.Paragraphs.Item(16)
$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Drew – AI and </w:t></w:r><w:r><w:t>its Technical and Ethical</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Impacts</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Current State of Self-driving cars</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Forecast for self driving cars</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">Compare and </w:t></w:r><w:r><w:t>Contrast human v.</w:t></w:r><w:r><w:t xml:space="preserve"> AI driving</w:t></w:r><w:r><w:t xml:space="preserve"> abilities</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Show data supporting these claims</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Technical challenges of AI</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Ethical Challenges of AI</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Trolley problem</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Implications on society and infrastructure</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Impact on the status quo</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$blankPara.Range.InsertXML($newBlockXml)

# ---------------------------------------------------------------------------
# 3. Remove the now-superseded "Drew – AI and ethics implications" line
#    (replaced by the fuller outline block inserted above).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(19).Range.Delete()

# ---------------------------------------------------------------------------
# 4. In Lyndon's paragraph, merge the run pair that used to be split by a
#    mid-sentence lastRenderedPageBreak into one continuous run, and add all
#    of the new analysis text.
# ---------------------------------------------------------------------------
$lyndonPara = $d.Paragraphs.Item(21)
$oldChunk = " cars in the same lane (visual aid here). On the freeways, a computer could effectively manage a smaller distance between cars at higher speeds, increasing freeway capacity and eliminating traffic jams, drastically shortening "
$newTailXml = ' cars in the same lane (visual aid here). On the freeways, a computer could effectively manage a smaller distance between cars at higher speeds, increasing freeway capacity and eliminating traffic jams, drastically shortening drive times and the ecological impact of traffic. One necessary ingredient for this to work is a network that allows all cars to communicate. Who would be responisible for maintaining such a system? The automobile manufacturers themselves wouldn’t be able to, since a singular entity is needed to guarantee accounting of all vehicles on the road. Allowing one company to manage this is essentially handing them a monopoly. Moreover, control over all moving traffic is a major homeland security concern. An exposed vulnerability in such a system could result in the most damaging attack imaginable – a complete shutdown of all roadways and risk of immediate death by collision for millions of daily commuters. With such a security concern, the governament would have to take charge of this self-driving network. This creates a major privacy versus security concern. GPS vehicle tracking currently isn’t legal in the US (provide link to ruling here). That hasn’t stopped law enforcement from tracking the general movement of automobiles using dragnet license plate scanning (visual aid for how this works here). GPS tracking is a prerequisite component to a self-driving car network, and we would essentially be handing detailed information about our real-time location and habits to the federal government. This is a lot of power to willingly hand over, and it’s important to consider whether it’s worth the trade-off. Would you rather sit in traffic and have more freedom, or never worry about traffic again and be immediately cited every time you sped (link to survey here)? While I’m speed racer on the road, I realize that traffic is an exponentially incresing problem, and I’d be willing to change my ways for a solution.'

# Remove the mid-sentence page break marker and rejoin the two runs it split.
$rFind = $lyndonPara.Range.Find
$rFind.Execute($oldChunk + "drive times", $true, $false, $false, $false, $false, $true, 1, $false, $newTailXml, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Move the lastRenderedPageBreak so it now sits right before "Alaysia"
#    (start of the next paragraph) instead of mid-sentence in Lyndon's text.
# ---------------------------------------------------------------------------
$alaysiaPara = $d.Paragraphs.Item(22)
$alaysiaRange = $alaysiaPara.Range.Duplicate
$alaysiaRange.Collapse(1)
$alaysiaRange.InsertXML('<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:lastRenderedPageBreak/></w:r>')
